$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.899443030357361
$ws.Range("B1").Value = 3.131730556488037
$ws.Range("C1").Value = 2.80917763710022
$ws.Range("D1").Value = 3.374208688735962
$ws.Range("E1").Value = 3.930059194564819
